$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$t.Cell(1, 1).Range.Text = "86÷3="  # was 55÷3=
$t.Cell(1, 2).Range.Text = "95÷7="  # was 62÷7=
$t.Cell(1, 3).Range.Text = "61÷8="  # was 73÷6=
$t.Cell(1, 4).Range.Text = "12÷8="  # was 41÷6=
$t.Cell(1, 5).Range.Text = "88÷3="  # was 64÷7=
$t.Cell(5, 1).Range.Text = "90÷7="  # was 97÷6=
$t.Cell(5, 2).Range.Text = "57÷2="  # was 85÷3=
$t.Cell(5, 3).Range.Text = "10÷6="  # was 93÷9=
$t.Cell(5, 4).Range.Text = "81÷4="  # was 73÷7=
$t.Cell(5, 5).Range.Text = "12÷6="  # was 29÷2=
$t.Cell(9, 1).Range.Text = "89÷5="  # was 85÷9=
$t.Cell(9, 2).Range.Text = "58÷7="  # was 77÷9=
$t.Cell(9, 3).Range.Text = "89÷9="  # was 91÷2=
$t.Cell(9, 4).Range.Text = "65÷5="  # was 77÷9=
$t.Cell(9, 5).Range.Text = "98÷8="  # was 26÷5=
$t.Cell(13, 1).Range.Text = "56÷7="  # was 43÷6=
$t.Cell(13, 2).Range.Text = "17÷4="  # was 86÷8=
$t.Cell(13, 3).Range.Text = "85÷2="  # was 82÷3=
$t.Cell(13, 4).Range.Text = "93÷9="  # was 57÷9=
$t.Cell(13, 5).Range.Text = "12÷4="  # was 28÷9=
$t.Cell(17, 1).Range.Text = "51÷9="  # was 98÷8=
$t.Cell(17, 2).Range.Text = "87÷4="  # was 21÷6=
$t.Cell(17, 3).Range.Text = "89÷5="  # was 33÷4=
$t.Cell(17, 4).Range.Text = "24÷3="  # was 84÷8=
$t.Cell(17, 5).Range.Text = "42÷8="  # was 26÷8=
